$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append four new data rows (190-193) to the bottom of the table
$newRows = @(
    @('XS0852986156','ISIN','Corp','KY','CN','ACTV','Financial','#N/A Field Not Applicable','#N/A Field Not Applicable','BBB+','Baa1','#N/A N/A'),
    @('XS1856799421','ISIN','Corp','VG','CN','ACTV','Consumer, Non-cyclical','#N/A Field Not Applicable','#N/A Field Not Applicable','BBB','Baa1','#N/A N/A'),
    @('XS1125272143','ISIN','Corp','VG','HK','ACTV','Consumer, Cyclical','#N/A Field Not Applicable','#N/A Field Not Applicable','A-','A3','#N/A N/A'),
    @('XS2008566197','ISIN','Corp','VG','CN','ACTV','Industrial','#N/A Field Not Applicable','#N/A Field Not Applicable','#N/A N/A','#N/A N/A','BBB+')
)

$startRow = 190
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Update the view: selection + scroll position
$ws.Range("M187").Select()
